$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "BEPEfCT": replace the placeholder "Boolean" row with a full
# per-industry-sector breakdown table, add a unit/italic header cell, and
# adjust page layout
# ---------------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("BEPEfCT")

# Replace the old single "Boolean" row with the full sector table
$sectors = @(
    @("agriculture and forestry 01T03", 1),
    @("coal mining 05", 0),
    @("oil and gas extraction 06", 0),
    @("other mining and quarrying 07T08", 0),
    @("food beverage and tobacco 10T12", 0),
    @("textiles apparel and leather 13T15", 0),
    @("wood products 16", 0),
    @("pulp paper and printing 17T18", 0),
    @("refined petroleum and coke 19", 0),
    @("chemicals 20", 0),
    @("rubber and plastic products 22", 0),
    @("glass and glass products 231", 0),
    @("cement and other nonmetallic minerals 239", 0),
    @("iron and steel 241", 0),
    @("other metals 242", 0),
    @("metal products except machinery and vehicles 25", 0),
    @("computers and electronics 26", 0),
    @("appliances and electrical equipment 27", 0),
    @("other machinery 28", 0),
    @("road vehicles 29", 0),
    @("nonroad vehicles 30", 0),
    @("other manufacturing 31T33", 0),
    @("energy pipelines and gas processing 352T353", 0),
    @("water and waste 36T39", 1),
    @("construction 41T43", 0)
)

$row = 2
foreach ($sector in $sectors) {
    $wsBool.Cells.Item($row, 1).Value = $sector[0]
    $wsBool.Cells.Item($row, 2).Value = $sector[1]
    $row = $row + 1
}

# New italic "Unit" header cell in A1 (B1 keeps its existing text)
$wsBool.Range("A1").Value = "Unit: boolean (0 or 1)"
$wsBool.Range("A1").Font.Italic = $true

# Widen column A so the sector labels are fully visible (closest the engine's
# pixel-snapped column-width model can get to the authored 47.1796875 value)
$wsBool.Columns.Item(1).ColumnWidth = 46.3

# Match the page orientation recorded in the saved file
$wsBool.PageSetup.Orientation = 1

# Leave the selection on B5, as in the saved workbook
$wsBool.Range("B5").Select()

# ---------------------------------------------------------------------------
# Sheet "About": add the two new explanatory rows at the bottom (13 & 14)
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A13").Value = 'In the U.S., we exempt agriculture and water and waste process emissions. Generally, '
$wsAbout.Range("A14").Value = 'proposed taxes do not cover these sectors.'

# Leave the selection where Excel would land after typing into A14 (next row)
$wsAbout.Range("A15").Select()
